$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 16).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 16)
    if ($cell.Value2 -eq "2022-06-21 06:53:08") {
        $cell.Value = "2022-06-21 07:37:36"
    }
}
